$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price / 1h-volume figures with the latest scraped values.
# Some prices are plain decimals (e.g. 85.20); prefixing with a leading apostrophe keeps
# them stored as text (matching the other price cells, e.g. "42.869.35") instead of
# letting Excel coerce them into numbers, which would silently drop trailing zeros.

$ws.Range("D2").Value = "42.869.35"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "2.538.94"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'315.37"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'96.18"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  -1.35%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").Value = "'36.08"
$ws.Range("E10").Value = "  -0.79%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("E13").Value = "  -2.62%  "
$ws.Range("D14").Value = "2.925.44"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "2.546.16"
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("D16").Value = "'15.23"
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("D17").Value = "'0.853"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "42.925.84"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").Value = "'13.07"
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("D20").Value = "'6.84"
$ws.Range("E20").Value = "  +3.72%  "
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "'70.10"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Value = "'253.23"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -1.17%  "
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("D26").Value = "'26.66"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +2.22%  "
$ws.Range("D29").Value = "'40.99"
$ws.Range("E29").Value = "  +7.72%  "
$ws.Range("D30").Value = "'10.44"
$ws.Range("E30").Value = "  +3.16%  "
$ws.Range("D31").Value = "'5.94"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "'157.46"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  +4.07%  "
$ws.Range("D34").Value = "'3.34"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "'19.02"
$ws.Range("E35").Value = "  -3.04%  "
$ws.Range("E36").Value = "  +2.22%  "
$ws.Range("D37").Value = "'0.0783"
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("E38").Value = "  -0.76%  "
$ws.Range("E39").Value = "  -1.12%  "
$ws.Range("D40").Value = "'2.32"
$ws.Range("E40").Value = "  +13.48%  "
$ws.Range("D41").Value = "'23.03"
$ws.Range("E41").Value = "  -4.47%  "
$ws.Range("D42").Value = "'3.86"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "'0.0304"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "2.034.44"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "'9.17"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("D48").Value = "'85.20"
$ws.Range("D49").Value = "'106.43"
$ws.Range("E49").Value = "  +4.83%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.780.03"
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "'74.68"
$ws.Range("E51").Value = "  +0.08%  "
